$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new data rows before the old row 19 ("item 15"), pushing it (and
# everything below) down. Each new row copies the formatting of row 18 and
# gets the next sequential item number in column A.
for ($i = 0; $i -lt 8; $i++) {
    $targetRow = 19 + $i
    $ws.Rows.Item($targetRow).Insert(-4121)
    $ws.Range("A18:K18").Copy()
    $ws.Range("A" + $targetRow + ":K" + $targetRow).PasteSpecial(-4122)
    $ws.Cells.Item($targetRow, 1).Value = 15 + $i
}
$excel.CutCopyMode = $false

# The old row 19 ("item 15") is now row 27 - renumber it to item 23, since it
# is now the 23rd (and last) data row.
$ws.Cells.Item(27, 1).Value = 23

# The print area needs to grow to cover the newly inserted rows too.
$ws.PageSetup.PrintArea = "`$B`$1:`$K`$31"

# Leave the selection on the new TOTAL row's "Expected Insolation" cell.
$ws.Range("E28").Select() | Out-Null

Write-Output "done"
